# Insert a new header row at the top of the sheet ("ENSEMBL_ID"),
# pushing the existing 369 ENSG identifiers down one row (A2:A370),
# matching the new HKG list comparison described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "ENSEMBL_ID"
